# Final draft text & metadata edits before review
# Applies the Personnel-sheet data corrections:
#  - Swap row 6 / row 7 contents (Danielle Aldrett <-> Kevin Cahill), moving
#    Kevin Cahill's e-mail/hyperlink from E6 to E7.
#  - Update the personnel entry in row 8 from "Justin Ossolinski" to
#    "S. Alejandra Casillo Cieza".
#  - Update the sheet selection to match the author's final cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Personnel")

# --- Row 6 / Row 7: swap Kevin Cahill and Danielle Aldrett ------------------
# Row 6 becomes Danielle Aldrett (no e-mail hyperlink in this row any more)
$ws.Range("A6").Value2 = "Danielle"
$ws.Range("C6").Value2 = "Aldrett"
$ws.Range("E6").Clear()

# Row 7 becomes Kevin Cahill, including his e-mail + hyperlink
$ws.Range("A7").Value2 = "Kevin"
$ws.Range("C7").Value2 = "Cahill"
$ws.Range("E7").Value2 = "kcahill@whoi.edu"
$ws.Range("E7").Style = "Hyperlink"

# --- Row 8: Justin Ossolinski -> S. Alejandra Casillo Cieza -----------------
$ws.Range("A8").Value2 = "S. Alejandra"
$ws.Range("C8").Value2 = "Casillo Cieza"

# --- Rebuild the hyperlinks collection so the mailto link follows Kevin -----
# Cahill's row (E6 -> E7). The collection has to be rebuilt wholesale because
# individual Hyperlink.Delete() calls are not observed on save.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("E3"), "mailto:rstanle2@wellesley.edu")
$ws.Hyperlinks.Add($ws.Range("E4"), "mailto:rstanle2@wellesley.edu")
$ws.Hyperlinks.Add($ws.Range("E5"), "mailto:rstanle2@wellesley.edu")
$ws.Hyperlinks.Add($ws.Range("E7"), "mailto:kcahill@whoi.edu")

# Re-assert the Hyperlink style on the linked cells: Hyperlinks.Add() above
# can otherwise register a redundant duplicate style for them.
$ws.Range("E3").Style = "Hyperlink"
$ws.Range("E4").Style = "Hyperlink"
$ws.Range("E5").Style = "Hyperlink"
$ws.Range("E7").Style = "Hyperlink"

# --- Selection / view state ---------------------------------------------
$ws.Range("E23").Select()

Write-Host "Personnel sheet updated"
